$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header ("Experiment Group" must become shared-string index 21,
# so it needs to be written before any of the other brand-new strings below)
$ws.Range("E1").Value = "Experiment Group"

# Mark existing scale-factor-0.5 rows (2-11) as Experiment Group 1
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Value = 1
}

# Helper to write a value as literal TEXT even when it looks numeric
# (e.g. "4.0"), without leaving the cell's style changed afterwards.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# New block: scale factor 4.0 results (Experiment Group 2), rows 12-21
$ws.Range("A12").Value = "Scale Factor"
$ws.Range("B12").Value = "BerlinMod Benchmark"
Set-TextValue $ws.Range("C12") "4.0"
Set-TextValue $ws.Range("D12") "4.0"
$ws.Range("E12").Value = 2

$ws.Range("A13").Value = "Total Size (GB)"
$ws.Range("B13").Value = "BerlinMod Benchmark"
$ws.Range("C13").Value = 32
$ws.Range("D13").Value = 32
$ws.Range("E13").Value = 2

$ws.Range("A14").Value = "Azure VM Name"
$ws.Range("B14").Value = "VM"
$ws.Range("C14").Value = "B4ms"
$ws.Range("D14").Value = "B4ms"
$ws.Range("E14").Value = 2

$ws.Range("A15").Value = "# of CPU(s)"
$ws.Range("B15").Value = "VM"
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = 2

$ws.Range("A16").Value = "RAM (GB)"
$ws.Range("B16").Value = "VM"
$ws.Range("C16").Value = 16
$ws.Range("D16").Value = 32
$ws.Range("E16").Value = 2

$ws.Range("A17").Value = "Disk Size (GB)"
$ws.Range("B17").Value = "VM"
$ws.Range("C17").Value = 60
$ws.Range("D17").Value = 120
$ws.Range("E17").Value = 2

$ws.Range("A18").Value = "Disk Type"
$ws.Range("B18").Value = "VM"
$ws.Range("C18").Value = "Premium SSD"
$ws.Range("D18").Value = "Premium SSD"
$ws.Range("E18").Value = 2

$ws.Range("A19").Value = "Concurrent Users"
$ws.Range("B19").Value = "Scalar Benchmark"
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 2

$ws.Range("A20").Value = "Total Duration (hours)"
$ws.Range("B20").Value = "Scalar Benchmark"
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 2

$ws.Range("A21").Value = "User Request"
$ws.Range("B21").Value = "Scalar Benchmark"
$ws.Range("C21").Value = "BerlinMod Query14"
$ws.Range("D21").Value = "BerlinMod Query14"
$ws.Range("E21").Value = 2

# Set column E width to match the bestFit applied to the new column
$ws.Columns.Item(5).ColumnWidth = 17.42578125

# Update selection to match the final state
$ws.Range("D21").Select()
